# Insert a new data row at row 27 (pushes existing rows 27-101 down to 28-102,
# shifting their values/formats along; the worksheet's used range grows from
# A1:R101 to A1:R102 automatically).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new record's values.
$ws.Range("A27").Value = 10
$ws.Range("B27").Value = 'Vega Modelo de Temuco'
$ws.Range("C27").Value = 'La Araucanía'
$ws.Range("D27").Value = 44414
$ws.Range("E27").Value = 9
$ws.Range("F27").Value = 100112005
$ws.Range("G27").Value = 'Puerro'
$ws.Range("H27").Value = 'Azul de Maquehue'
$ws.Range("I27").Value = 'Primera'
$ws.Range("J27").Value = 40
$ws.Range("K27").Value = 8000
$ws.Range("L27").Value = 8000
$ws.Range("M27").Value = 8000
$ws.Range("N27").Value = '$/docena de paquetes'
$ws.Range("O27").Value = 'Provincia de Cautín'
$ws.Range("P27").Value = 667
$ws.Range("Q27").Value = 12
$ws.Range("R27").Value = 'Hortaliza'
